# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement table (rows 16-30) listed 15 overdue periods
# (2312, 2401..2502). It is replaced with 13 periods, newest first,
# ending on the oldest one (2412, 2411, ..., 2401, 2312), and the
# summary figures (overdue value / period count) are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last data row (row 30) carries the special "bottom of table" border
# formatting. Grab a copy of it now, before any row is touched, so we can
# stamp it onto the new last row (row 28) once the two trailing rows are
# removed.
$ws.Range("B30:J30").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)

# New period list for rows 16-28 (newest -> oldest), replacing the old
# 2312, 2401, 2402, ..., 2412, 2501, 2502.
$periods = @("2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401","2312")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Drop the two rows that fell off the end of the new, shorter period list
# (used to hold periods 2501 and 2502). This shifts everything below
# (the blank spacer + signature block) up by two rows, e.g. old rows
# 35/36 become 33/34.
$ws.Rows("29:30").Delete()

# Update the summary figures to match the new data set.
$ws.Range("E11").Value = 832000
$ws.Range("F13").Value = 13
